$d = $word.ActiveDocument

# 1) Flag every drawing's run as "no proofing" (<w:rPr><w:noProof/></w:rPr>)
#    - inline pictures live in InlineShapes
foreach ($shp in $d.InlineShapes) {
    $shp.Range.NoProofing = $true
}

#    - floating (anchored) pictures don't resolve a usable Range through
#      Shape.Anchor in this host, so flag every paragraph whose XML proves
#      it actually holds a <w:drawing> run (covers the two anchored photos
#      that share one paragraph, plus is a no-op anywhere else).
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs($i)
    if ($para.Range.WordOpenXML -like "*<w:drawing*") {
        $para.Range.NoProofing = $true
    }
}

# 2) Wording tweak: the queue is "paused", not "stopped"
$d.Content.Find.Execute("stopped", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "paused", 2)
